$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 33   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/12/2026  Through  1/18/2026"

# --- Row 15 ---
$ws.Range("C15").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("F15").Value = "0"
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -100
$ws.Range("I15").Value = "0"
$ws.Range("J15").Value = 3

# --- Row 16 ---
$ws.Range("F16").Value = 3
$ws.Range("H16").Value = -70
$ws.Range("I16").Value = 2
$ws.Range("J16").Value = 5
$ws.Range("K16").Value = -60
$ws.Range("L16").Value = -50
$ws.Range("M16").Value = -71.428571428571
$ws.Range("N16").Value = -92.857142857142

# --- Row 17 ---
$ws.Range("C17").Value = 4
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 21
$ws.Range("H17").Value = 61.538461538461
$ws.Range("I17").Value = 12
$ws.Range("J17").Value = 9
$ws.Range("K17").Value = 33.333333333333
$ws.Range("L17").Value = 71.428571428571
$ws.Range("M17").Value = 20
$ws.Range("N17").Value = 71.428571428571

# --- Row 18 ---
$ws.Range("D16").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = 2
$ws.Range("D16").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = 1
$ws.Range("E16").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 2
$ws.Range("H18").Value = 100
$ws.Range("D16").Copy()
$ws.Range("I18").PasteSpecial(-4122)
$ws.Range("I18").Value = 2
$ws.Range("D16").Copy()
$ws.Range("J18").PasteSpecial(-4122)
$ws.Range("J18").Value = 1
$ws.Range("E16").Copy()
$ws.Range("K18").PasteSpecial(-4122)
$ws.Range("K18").Value = 100
$ws.Range("L18").Value = -33.333333333333
$ws.Range("M18").Value = -86.666666666666
$ws.Range("N18").Value = -94.117647058823

# --- Row 19 ---
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 15
$ws.Range("G19").Value = 11
$ws.Range("H19").Value = 36.363636363636
$ws.Range("I19").Value = 7
$ws.Range("J19").Value = 8
$ws.Range("K19").Value = -12.5
$ws.Range("L19").Value = -65
$ws.Range("M19").Value = -12.5
$ws.Range("N19").Value = -66.666666666666

# --- Row 20 ---
$ws.Range("C20").Value = 2
$ws.Range("C15").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value = "0"
$ws.Range("L15").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = "***.*"
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -12.5
$ws.Range("I20").Value = 5
$ws.Range("K20").Value = 25
$ws.Range("L20").Value = -37.5
$ws.Range("M20").Value = -50
$ws.Range("N20").Value = -95.652173913043

# --- Row 21 ---
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = 20
$ws.Range("F21").Value = 48
$ws.Range("G21").Value = 48
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 28
$ws.Range("J21").Value = 30
$ws.Range("K21").Value = -6.666666666666
$ws.Range("L21").Value = -34.883720930232
$ws.Range("M21").Value = -44
$ws.Range("N21").Value = -86.473429951690

# --- Row 22 ---
$ws.Range("C15").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("G22").Value = "0"
$ws.Range("L15").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("H22").Value = "***.*"
$ws.Range("E16").Copy()
$ws.Range("M22").PasteSpecial(-4122)
$ws.Range("M22").Value = -100

# --- Row 23 ---
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 5
$ws.Range("H23").Value = -37.5
$ws.Range("I23").Value = 4
$ws.Range("J23").Value = 4
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = -60
$ws.Range("M23").Value = 300

# --- Row 24 ---
$ws.Range("C24").Value = 7
$ws.Range("D24").Value = 6
$ws.Range("E24").Value = 16.666666666666
$ws.Range("F24").Value = 40
$ws.Range("G24").Value = 30
$ws.Range("H24").Value = 33.333333333333
$ws.Range("I24").Value = 28
$ws.Range("J24").Value = 16
$ws.Range("K24").Value = 75
$ws.Range("L24").Value = -6.666666666666
$ws.Range("M24").Value = 16.666666666666

# --- Row 25 ---
$ws.Range("C15").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("C25").Value = "0"
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = -100
$ws.Range("F25").Value = 6
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = -50
$ws.Range("J25").Value = 9
$ws.Range("K25").Value = -55.555555555555
$ws.Range("L25").Value = -50

# --- Row 26 ---
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 75
$ws.Range("F26").Value = 26
$ws.Range("G26").Value = 23
$ws.Range("H26").Value = 13.043478260869
$ws.Range("I26").Value = 17
$ws.Range("J26").Value = 15
$ws.Range("K26").Value = 13.333333333333
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -22.727272727272

# --- Row 27 ---
$ws.Range("C15").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("F27").Value = "0"
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -100
$ws.Range("I27").Value = "0"
$ws.Range("J27").Value = 4

# --- Row 28 ---
$ws.Range("D16").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 2
$ws.Range("E16").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = 50
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 133.333333333333
$ws.Range("I28").Value = 7
$ws.Range("J28").Value = 3
$ws.Range("K28").Value = 133.333333333333
$ws.Range("L28").Value = 75

# --- Row 29 ---
$ws.Range("C15").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value = "0"
$ws.Range("L15").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = "***.*"

# --- Row 30 ---
$ws.Range("C15").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = "0"
$ws.Range("L15").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = "***.*"

# --- Row 33 ---
$ws.Range("C15").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("D33").Value = "0"
$ws.Range("L15").Copy()
$ws.Range("E33").PasteSpecial(-4122)
$ws.Range("E33").Value = "***.*"
